$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds plain-text price strings; force text format while writing
# numeric-looking values so Excel does not auto-convert them to numbers,
# then restore the default (unstyled) cell style.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.602.37"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").Value = "1.821.19"
$ws.Range("E3").Value = "  +1.59%  "
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("D6").Value = "305.13"
$ws.Range("E6").Value = "  -0.70%  "
$ws.Range("D7").Value = "0.4664"
$ws.Range("E7").Value = "  +2.29%  "
$ws.Range("D8").Value = "0.3589"
$ws.Range("E8").Value = "  -0.88%  "
$ws.Range("D9").Value = "0.07116"
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("D10").Value = "0.8992"
$ws.Range("E10").Value = "  +2.35%  "
$ws.Range("D11").Value = "0.07797"
$ws.Range("E11").Value = "  -0.48%  "
$ws.Range("D12").Value = "19.35"
$ws.Range("E12").Value = "  -0.89%  "
$ws.Range("D13").Value = "1.850.28"
$ws.Range("E13").Value = "  +3.25%  "
$ws.Range("D14").Value = "5.250"
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("D15").Value = "6.331"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("E16").Value = "  +2.58%  "
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("D18").Value = "0.000008541"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("D20").Value = "26.643.76"
$ws.Range("E20").Value = "  +0.76%  "
$ws.Range("E21").Value = "  -0.94%  "
$ws.Range("E22").Value = "  +0.35%  "
$ws.Range("E23").Value = "  +0.48%  "
$ws.Range("D24").Value = "1.932"
$ws.Range("E24").Value = "  -2.00%  "
$ws.Range("D25").Value = "152.00"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").Value = "1.966"
$ws.Range("E27").Value = "  -3.56%  "
$ws.Range("D28").Value = "113.66"
$ws.Range("E28").Value = "  +1.49%  "
$ws.Range("E29").Value = "  -1.33%  "
$ws.Range("D30").Value = "0.08785"
$ws.Range("E30").Value = "  +1.51%  "
$ws.Range("D31").Value = "3.118"
$ws.Range("E31").Value = "  +1.56%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "0.7272"
$ws.Range("E32").Value = "  +0.55%  "
$ws.Range("B33").Value = "RenderToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D33").Value = "2.723"
$ws.Range("E33").Value = "  +2.27%  "
$ws.Range("D34").Value = "4.424"
$ws.Range("E34").Value = "  -0.48%  "
$ws.Range("D35").Value = "1.120"
$ws.Range("E35").Value = "  +0.98%  "
$ws.Range("D36").Value = "1.075"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("D37").Value = "0.01926"
$ws.Range("E37").Value = "  -0.96%  "
$ws.Range("D38").Value = "2.918"
$ws.Range("E38").Value = "  +1.58%  "
$ws.Range("D39").Value = "0.05095"
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("D40").Value = "6.810"
$ws.Range("E40").Value = "  -1.19%  "
$ws.Range("D41").Value = "0.5021"
$ws.Range("E41").Value = "  -4.23%  "
$ws.Range("D42").Value = "0.1489"
$ws.Range("E42").Value = "  -2.50%  "
$ws.Range("D43").Value = "7.949"
$ws.Range("E43").Value = "  -0.80%  "
$ws.Range("D44").Value = "1.009"
$ws.Range("E44").Value = "  +0.32%  "
$ws.Range("D45").Value = "0.4632"
$ws.Range("E45").Value = "  -1.29%  "
$ws.Range("D46").Value = "9.971"
$ws.Range("E46").Value = "  +0.87%  "
$ws.Range("D47").Value = "98.24"
$ws.Range("E47").Value = "  -1.77%  "
$ws.Range("E48").Value = "  -2.20%  "
$ws.Range("D49").Value = "0.05994"
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("D50").Value = "63.56"
$ws.Range("E50").Value = "  -0.76%  "
$ws.Range("D51").Value = "35.64"
$ws.Range("E51").Value = "  -1.68%  "

$ws.Range("D2:D51").Style = "Normal"
